$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label for the Disk Utility column
$ws.Range("O2").Value = "Utility (Percent)"

# Append " msec" unit suffix to Read Latency values (min/max/average) for data rows 3-23
for ($r = 3; $r -le 23; $r++) {
    foreach ($col in @("I", "J", "K")) {
        $cellRef = "$col$r"
        $current = $ws.Range($cellRef).Text
        $ws.Range($cellRef).Value = "$current msec"
    }
}
